# "Visor indicadores en tableau"
# - Adds new Tableau-sourced metrics (extension, creation year, zone, notes)
#   for the last three London boroughs on the "Barrios" sheet (rows 31-33).
# - Switches the active tab from "Barrios" back to "Sheet1" and updates the
#   selection/scroll position left behind on "Barrios".

$wb = $excel.ActiveWorkbook

$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsBarrios = $wb.Worksheets.Item("Barrios")

# --- New data for rows 31-33 on "Barrios" ---

# Row 31 (Havering)
$wsBarrios.Range("C31").Value = 38.82
$wsBarrios.Range("D31").Value = 1965
$wsBarrios.Range("E31").Value = "Outer London"
$wsBarrios.Range("F31").Value = "Sede de las olimpiadas con Lee Valley Hockey and Tennis Centre and part of the Queen Elizabeth Olympic Park"

# Row 32 (Hillingdon)
$wsBarrios.Range("C32").Value = 34.26
$wsBarrios.Range("E32").Value = "Inner London"
$wsBarrios.Range("F32").Value = "Sede del nuevo mercado de Covent Garden y el Helipuerto de Londres"

# Row 33 (Hounslow)
$wsBarrios.Range("C33").Value = 21.48
$wsBarrios.Range("E33").Value = "Inner London"
$wsBarrios.Range("F33").Value = "Lugares emblemáticos como el Parlamento, El British government, zona de compras Oxford Street, Regent Street, Picacadilly, Bond Street, Soho, Buckingham Palace, Westminster Abbey, WhiteHall, Trafalgar Square, Hyde Park, gran parte de Regent's Park"

# --- View / selection changes ---

# Leave "Barrios" scrolled to the top with D8 selected (was topLeftCell A22 / C31 selected)
$wsBarrios.Range("D8").Select()

# Make "Sheet1" the active tab again (was "Barrios")
$wsSheet1.Activate()
